# Append a new data row (row 35) to each of the four worksheets,
# mirroring the structure of the existing rows (A: timestamp w/ style,
# B-E: inline hex strings, F-I: numeric decimals).

$wb = $excel.ActiveWorkbook

$rowsData = @{
    "MID_LFT_#1" = @{
        A = 45821.46068287037
        B = "0x01,0x90"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1b,0x41,0x0c,"
        D = "0x01,0x80"
        E = "0x07"
        F = 400
        G = [double]"5.68631262647113e+23"
        H = 384
        I = 7
    }
    "MID_LFT_#2" = @{
        A = 45821.46068287037
        B = "0x01,0x7c"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
        D = "0x01,0x74"
        E = "0x19"
        F = 380
        G = [double]"5.68432987514711e+23"
        H = 372
        I = 25
    }
    "MID_PLT_#1" = @{
        A = 45821.46068287037
        B = "0x00,0x6e"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c,"
        D = "0x00,0x6C"
        E = "0x15"
        F = 110
        G = [double]"5.68631262647113e+23"
        H = 108
        I = 15
    }
    "MID_PLT_#2" = @{
        A = 45821.46068287037
        B = "0x00,0x82"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c,"
        D = "0x00,0x80"
        E = "0x9"
        F = 130
        G = [double]"5.68631262647113e+23"
        H = 128
        I = 9
    }
}

foreach ($ws in $wb.Worksheets) {
    $data = $rowsData[$ws.Name]
    if ($data -eq $null) { continue }

    $newRow = 35

    $ws.Cells.Item($newRow, 1).Value = $data.A
    $ws.Cells.Item($newRow, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Cells.Item($newRow, 2).Value = $data.B
    $ws.Cells.Item($newRow, 3).Value = $data.C
    $ws.Cells.Item($newRow, 4).Value = $data.D
    $ws.Cells.Item($newRow, 5).Value = $data.E

    $ws.Cells.Item($newRow, 6).Value = $data.F
    $ws.Cells.Item($newRow, 7).Value = $data.G
    $ws.Cells.Item($newRow, 8).Value = $data.H
    $ws.Cells.Item($newRow, 9).Value = $data.I
}
